# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" for the first localized file
# (cc14e37b-...) row in each per-language status sheet, reflecting a
# newly-generated handoff report.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("D2").Value = "2016-02-22 10:48:15"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("D2").Value = "2016-02-22 10:48:28"
